$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shift existing D:K to F:M)
$ws.Columns("D:E").Insert()

# Copy number formats from the (now-shifted) original D:E columns (now F:G) into the new D:E columns
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the two new quarter columns (Dec-2018 / Sep-2018) with their reported values
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 97900
$ws.Range("E8").Value = 97000
$ws.Range("D9").Value = 63700
$ws.Range("E9").Value = 62300
$ws.Range("D10").Value = 34200
$ws.Range("E10").Value = 34700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 1700
$ws.Range("E15").Value = 900
$ws.Range("D17").Value = 83500
$ws.Range("E17").Value = 81700
$ws.Range("D18").Value = 14400
$ws.Range("E18").Value = 15300
$ws.Range("D20").Value = -13300
$ws.Range("E20").Value = -1700
$ws.Range("D21").Value = 6700
$ws.Range("E21").Value = 19600
$ws.Range("D22").Value = 5600
$ws.Range("E22").Value = 5100
$ws.Range("D23").Value = -4600
$ws.Range("E23").Value = 8400
$ws.Range("D24").Value = -200
$ws.Range("E24").Value = 2300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -4400
$ws.Range("E26").Value = 6200
$ws.Range("D27").Value = -4300
$ws.Range("E27").Value = 6300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 13300
$ws.Range("E32").Value = 1700
$ws.Range("D33").Value = -4300
$ws.Range("E33").Value = 6300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -4300
$ws.Range("E35").Value = 6300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 33000
$ws.Range("E41").Value = 28000
$ws.Range("D42").Value = 1200
$ws.Range("E42").Value = 1500
$ws.Range("D43").Value = 147000
$ws.Range("E43").Value = 145700
$ws.Range("D44").Value = 91800
$ws.Range("E44").Value = 88500
$ws.Range("D45").Value = 20300
$ws.Range("E45").Value = 21400
$ws.Range("D46").Value = 293400
$ws.Range("E46").Value = 285100
$ws.Range("D47").Value = 7000
$ws.Range("E47").Value = 5500
$ws.Range("D48").Value = 149200
$ws.Range("E48").Value = 163500
$ws.Range("D49").Value = 32600
$ws.Range("E49").Value = 33400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 7600
$ws.Range("E52").Value = 3100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 489800
$ws.Range("E54").Value = 490600
$ws.Range("D57").Value = 65500
$ws.Range("E57").Value = 62500
$ws.Range("D58").Value = 21600
$ws.Range("E58").Value = 16100
$ws.Range("D59").Value = 35500
$ws.Range("E59").Value = 31600
$ws.Range("D60").Value = 122600
$ws.Range("E60").Value = 110200
$ws.Range("D61").Value = 229800
$ws.Range("E61").Value = 228400
$ws.Range("D62").Value = 4100
$ws.Range("E62").Value = 4700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 357400
$ws.Range("E66").Value = 344300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 11800
$ws.Range("E72").Value = 21400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 132400
$ws.Range("E76").Value = 146300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -4300
$ws.Range("E81").Value = 6300
$ws.Range("D83").Value = 5700
$ws.Range("E83").Value = 6000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 4400
$ws.Range("E89").Value = -3000
$ws.Range("D91").Value = -5900
$ws.Range("E91").Value = -2300
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -5800
$ws.Range("E94").Value = -1700
$ws.Range("D96").Value = -700
$ws.Range("E96").Value = -700
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 8100
$ws.Range("E100").Value = 3200
$ws.Range("D101").Value = -1600
$ws.Range("E101").Value = -400
$ws.Range("D102").Value = 5100
$ws.Range("E102").Value = -2000

# Data correction carried over from the source feed for Capital Expenditures (row 91)
$ws.Range("J91").Value = -2300
